$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TestCases")

# ---------------------------------------------------------------------------
# Row 15 ("wlanStatus with callback after disabling wlan profile" / G15):
# capture its current text BEFORE any other edits - this text is moving to
# the "Change Start page" row (row 2) slot is NOT what happens; rather this
# text is the one that was previously the very-last shared string and is now
# reused as-is for row 15 (its value does not change, only its backing
# shared-string slot shifts around as other strings change). We simply must
# make sure the final text in G15 equals this original text.
# ---------------------------------------------------------------------------
$g15Text = $ws.Range("G15").Value2

# ---------------------------------------------------------------------------
# Row 2 ("Change Start page"): G2's old "SetStartPage" snippet is replaced by
# a new, longer config-driven snippet (this is the "new way to change
# config" from the commit message).
# ---------------------------------------------------------------------------
$newG2 = @'
wait(3);
PullConfigxml;
ChangeConfigxml(Configuration/Applications/Application/General,StartPage,<StartPage value="http://127.0.0.1:8082/app/" name="Menu"/>);
ChangeConfigxml(Configuration,WebServer,<WebServer>);
ChangeConfigxml(Configuration/WebServer,Enabled,<Enabled VALUE="1"/>);
ChangeConfigxml(Configuration/WebServer,Port,<Port VALUE="8082"/>);
ChangeConfigxml(Configuration/WebServer,WebFolder,<WebFolder VALUE="\\auto\\ComplianceTest_JS\"/>);
ChangeConfigxml(Configuration/WebServer,Public,<Public VALUE="1"/>);
ChangeConfigxml(Configuration/Screen,FullScreen,<FullScreen value="0"/>);
PushConfigxml;
'@
$ws.Range("G2").Value = $newG2

# ---------------------------------------------------------------------------
# Row 4 ("wlanStatus with callback"): signalCallbackcount(...) call gains a
# second argument.
# ---------------------------------------------------------------------------
$g4 = $ws.Range("G4").Value2
$g4 = $g4.Replace("signalCallbackcount(results_id);", "signalCallbackcount(results_id,com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);")
$ws.Range("G4").Value = $g4

# ---------------------------------------------------------------------------
# Row 14 ("stopWlanStatus after calling wlanStatus"): checkstopwlanStatus(...)
# call gains a second argument.
# ---------------------------------------------------------------------------
$g14 = $ws.Range("G14").Value2
$g14 = $g14.Replace("checkstopwlanStatus(results_id);", "checkstopwlanStatus(results_id,com.symbol.enterprisebrowser/com.rhomobile.rhodes.RhodesActivity);")
$ws.Range("G14").Value = $g14

# ---------------------------------------------------------------------------
# Row 15: write back its original, unmodified text (content identical - only
# its shared-string slot changes implicitly on save).
# ---------------------------------------------------------------------------
$ws.Range("G15").Value = $g15Text

# ---------------------------------------------------------------------------
# Row heights grew because of the longer multi-line content.
# ---------------------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 332.25
$ws.Rows.Item(4).RowHeight = 315.75

# ---------------------------------------------------------------------------
# Active view: scrolled/selected down to G14 while editing.
# ---------------------------------------------------------------------------
$ws.Range("G14").Select()
$ws.Application.ActiveWindow.ScrollRow = 14
